$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# --- Sheet "VENTAS POR GRUPO" ---
# Zero out stale current-month figures that rolled off the monthly window
$ws1.Range("L2").Value = 0
$ws1.Range("C3").Value = 0
$ws1.Range("C4").Value = "0 de 2"
$ws1.Range("L4").Value = "0 de 2"

# --- Sheet "VENTA MENSUAL" ---
# Shift the rolling month headers forward by one month
$ws2.Range("C1").Value = "marzo"
$ws2.Range("D1").Value = "abril"
$ws2.Range("E1").Value = "mayo"
$ws2.Range("F1").Value = "junio"

# Shift the monthly figures forward to line up with the new headers
$ws2.Range("D2").Value = 144.53
$ws2.Range("E2").Value = 11.52
$ws2.Range("F2").Value = 0

$ws2.Range("E3").Value = 178.33
$ws2.Range("F3").Value = 0

$ws2.Range("D4").Value = 144.53
$ws2.Range("E4").Value = 189.85
$ws2.Range("F4").Value = 0

# Column widths follow the shifted columns (raw OOXML width = ColumnWidth + 0.83)
$ws2.Columns.Item(3).ColumnWidth = 10.17
$ws2.Columns.Item(4).ColumnWidth = 11.17
$ws2.Columns.Item(6).ColumnWidth = 10.17
